# Auto-generated script applying scheduled market-data refresh to the FFXIV Leve profit workbook.
# For each changed cell: numeric cells get their new value; cells that the refresh
# removed (no longer applicable, e.g. blank HQ profit because HQ price is 0) are cleared.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 162.5
$ws.Range("J43").Value = 164.28572
$ws.Range("L43").Value = 164.28572
$ws.Range("N43").Value = -302.28572
$ws.Range("H62").Value = 9146.538
$ws.Range("I62").Value = 7960
$ws.Range("J62").Value = 9888.125
$ws.Range("K62").Value = 7960
$ws.Range("L62").Value = 9888.125
$ws.Range("M62").Value = -7336
$ws.Range("N62").Value = -11136.125
$ws.Range("H65").Value = 9146.538
$ws.Range("I65").Value = 7960
$ws.Range("J65").Value = 9888.125
$ws.Range("K65").Value = 39800
$ws.Range("L65").Value = 49440.625
$ws.Range("M65").Value = -36680
$ws.Range("N65").Value = -55680.625
$ws.Range("H116").Value = 8205.77
$ws.Range("I116").Value = 3815.8333
$ws.Range("J116").Value = 11968.571
$ws.Range("K116").Value = 3815.8333
$ws.Range("L116").Value = 11968.571
$ws.Range("M116").Value = -373.8332999999998
$ws.Range("N116").Value = -18852.571
$ws.Range("H125").Value = 2966.6667
$ws.Range("J125").Value = 2966.6667
$ws.Range("L125").Value = 26700.0003
$ws.Range("N125").Value = -31620.0003
$ws.Range("H135").Value = 1441.5834
$ws.Range("I135").Value = 918.1667
$ws.Range("J135").Value = 1965
$ws.Range("K135").Value = 8263.5003
$ws.Range("L135").Value = 17685
$ws.Range("M135").Value = -5728.5003
$ws.Range("N135").Value = -22755
$ws.Range("H137").Value = 3009.0667
$ws.Range("J137").Value = 1556
$ws.Range("L137").Value = 4668
$ws.Range("N137").Value = -9768
$ws.Range("H141").Value = 9503.706
$ws.Range("I141").Value = 3320.75
$ws.Range("K141").Value = 9962.25
$ws.Range("M141").Value = -4782.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 999.5
$ws.Range("I2").Value = 999.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 999.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -886.5
$ws.Range("N2").ClearContents()
$ws.Range("H6").Value = 2000
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H23").Value = 20000
$ws.Range("J23").Value = 20000
$ws.Range("L23").Value = 20000
$ws.Range("N23").Value = -20518
$ws.Range("H61").Value = 3635.8823
$ws.Range("I61").Value = 5052.8887
$ws.Range("J61").Value = 2041.75
$ws.Range("K61").Value = 5052.8887
$ws.Range("L61").Value = 2041.75
$ws.Range("M61").Value = -4840.8887
$ws.Range("N61").Value = -2465.75
$ws.Range("H110").Value = 2410
$ws.Range("I110").Value = 2410
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 2410
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -365
$ws.Range("N110").ClearContents()
$ws.Range("H116").Value = 999.5
$ws.Range("I116").Value = 999.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 999.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1294.5
$ws.Range("N116").ClearContents()
$ws.Range("H136").Value = 3635.8823
$ws.Range("I136").Value = 5052.8887
$ws.Range("J136").Value = 2041.75
$ws.Range("K136").Value = 15158.6661
$ws.Range("L136").Value = 6125.25
$ws.Range("M136").Value = -12608.6661
$ws.Range("N136").Value = -11225.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 999.5
$ws.Range("I3").Value = 999.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 999.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -885.5
$ws.Range("N3").ClearContents()
$ws.Range("H22").Value = 1737.091
$ws.Range("I22").Value = 1889.5555
$ws.Range("J22").Value = 1051
$ws.Range("K22").Value = 1889.5555
$ws.Range("L22").Value = 1051
$ws.Range("M22").Value = -1716.5555
$ws.Range("N22").Value = -1397
$ws.Range("H134").Value = 3149.092
$ws.Range("I134").Value = 1357.86
$ws.Range("J134").Value = 6593.769
$ws.Range("K134").Value = 4073.58
$ws.Range("L134").Value = 19781.307
$ws.Range("M134").Value = -1538.58
$ws.Range("N134").Value = -24851.307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 650.13336
$ws.Range("I22").Value = 650.2
$ws.Range("J22").Value = 650
$ws.Range("K22").Value = 650.2
$ws.Range("L22").Value = 650
$ws.Range("M22").Value = -300.2
$ws.Range("N22").Value = -1350
$ws.Range("H58").Value = 5230.3706
$ws.Range("I58").Value = 2423.2942
$ws.Range("J58").Value = 10002.4
$ws.Range("K58").Value = 2423.2942
$ws.Range("L58").Value = 10002.4
$ws.Range("M58").Value = -2220.2942
$ws.Range("N58").Value = -10408.4
$ws.Range("H132").Value = 2445.389
$ws.Range("I132").Value = 1530.7059
$ws.Range("J132").Value = 3263.7896
$ws.Range("K132").Value = 4592.1177
$ws.Range("L132").Value = 9791.3688
$ws.Range("M132").Value = -2062.1177
$ws.Range("N132").Value = -14851.3688
$ws.Range("H136").Value = 5230.3706
$ws.Range("I136").Value = 2423.2942
$ws.Range("J136").Value = 10002.4
$ws.Range("K136").Value = 7269.882599999999
$ws.Range("L136").Value = 30007.2
$ws.Range("M136").Value = -4719.882599999999
$ws.Range("N136").Value = -35107.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1304578.2
$ws.Range("I132").Value = 2779657.2
$ws.Range("K132").Value = 8338971.600000001
$ws.Range("M132").Value = -8336441.600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 696
$ws.Range("I22").Value = 759.6
$ws.Range("J22").Value = 671.53845
$ws.Range("K22").Value = 759.6
$ws.Range("L22").Value = 671.53845
$ws.Range("M22").Value = -464.6
$ws.Range("N22").Value = -1261.53845
$ws.Range("H27").Value = 696
$ws.Range("I27").Value = 759.6
$ws.Range("J27").Value = 671.53845
$ws.Range("K27").Value = 759.6
$ws.Range("L27").Value = 671.53845
$ws.Range("M27").Value = -652.6
$ws.Range("N27").Value = -885.53845
$ws.Range("H46").Value = 667275.3
$ws.Range("I46").Value = 523.3333
$ws.Range("J46").Value = 1111776.6
$ws.Range("K46").Value = 523.3333
$ws.Range("L46").Value = 1111776.6
$ws.Range("M46").Value = -335.3333
$ws.Range("N46").Value = -1112152.6
$ws.Range("H55").Value = 45671.953
$ws.Range("I55").Value = 171.4375
$ws.Range("J55").Value = 167006.67
$ws.Range("K55").Value = 171.4375
$ws.Range("L55").Value = 167006.67
$ws.Range("M55").Value = 1.5625
$ws.Range("N55").Value = -167352.67
$ws.Range("H132").Value = 51576.617
$ws.Range("I132").Value = 87142.664
$ws.Range("J132").Value = 4155.222
$ws.Range("K132").Value = 261427.992
$ws.Range("L132").Value = 12465.666
$ws.Range("M132").Value = -258897.992
$ws.Range("N132").Value = -17525.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H111").Value = 33644
$ws.Range("J111").Value = 33644
$ws.Range("L111").Value = 33644
$ws.Range("N111").Value = -41824

